# This workbook holds weekly price-report rows (rows 2-30) for
# "Hortaliza, Agrícola del Norte S.A. de Arica - Perejil".
# The update re-shuffles the per-record fields (Fecha, Volumen,
# Precio minimo/maximo/promedio, Unidad de comercializacion, Origen,
# Precio $/Kg, Kg o Unidades -> columns D and J:Q) across the existing
# rows, while columns A,B,C,E,F,G,H,I,R (which already hold identical
# values on every row) are left untouched.
#
# Mapping: new row -> source row (value to move into that row)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 22
    3  = 7
    4  = 28
    5  = 2
    6  = 23
    7  = 12
    8  = 27
    9  = 15
    10 = 18
    11 = 13
    12 = 26
    13 = 8
    14 = 6
    15 = 4
    16 = 17
    17 = 3
    18 = 10
    19 = 21
    20 = 24
    21 = 14
    22 = 19
    23 = 11
    24 = 5
    25 = 29
    26 = 25
    27 = 16
    28 = 30
    29 = 9
    30 = 20
}

# First, snapshot the full D:Q block (Fecha ... Kg o Unidades) for every
# data row, since the mapping is a permutation with no fixed points and
# writes must not clobber data that is still needed as a source.
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $snapshot[$row] = $ws.Range("D$row`:Q$row").Value2
}

# Now write every row's new contents from the snapshot.
foreach ($row in $mapping.Keys) {
    $srcRow = $mapping[$row]
    $ws.Range("D$row`:Q$row").Value2 = $snapshot[$srcRow]
}
